$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Clear the (now pointless) highlight-fill style that MySQL!E3 had
#    - it was applied and then reset back to no fill, so drop the
#    explicit style and let the cell fall back to the default style.
# ------------------------------------------------------------------
$mysqlSheet = $wb.Worksheets.Item("MySQL")
$mysqlSheet.Range("E3").Style = "Normal"

# ------------------------------------------------------------------
# 2. Add the new "Sheet1" worksheet (comparative DB-systems summary)
#    after the existing "Cassandra" tab - it becomes the active tab.
# ------------------------------------------------------------------
$cassandraSheet = $wb.Worksheets.Item("Cassandra")
$newSheet = $wb.Worksheets.Add($null, $cassandraSheet)
$newSheet.Name = "Sheet1"

# Header row
$newSheet.Range("A1").Value = "System"
$newSheet.Range("B1").Value = "Read/Write"
$newSheet.Range("C1").Value = "Latency/Durability"
$newSheet.Range("D1").Value = "Sync/Async"
$newSheet.Range("E1").Value = "Row/Column optimized"
$newSheet.Range("F1").Value = "Replication"

# Cassandra row
$newSheet.Range("A2").Value = "Cassandra"
$newSheet.Range("B2").Value = "Write"
$newSheet.Range("C2").Value = "Tunable"
$newSheet.Range("D2").Value = "Tunable"
$newSheet.Range("E2").Value = "Column"
$newSheet.Range("F2").Value = "-"

# MySQL row
$newSheet.Range("A3").Value = "MySQL"
$newSheet.Range("B3").Value = "Read/Write"
$newSheet.Range("C3").Value = "Latency/Durability"
$newSheet.Range("D3").Value = "Sync"
$newSheet.Range("E3").Value = "Row"
$newSheet.Range("F3").Value = "Master-Slave"

# MongoDB row
$newSheet.Range("A4").Value = "MongoDB"
$newSheet.Range("B4").Value = "Read/Write"
$newSheet.Range("C4").Value = "Durability"
$newSheet.Range("D4").Value = "Sync/Async"
$newSheet.Range("E4").Value = "Document (Row)"
$newSheet.Range("F4").Value = "Replica Set"

# Fit the columns to their content and select the populated range,
# matching the sheet view Excel leaves behind after typing this table.
for ($col = 1; $col -le 6; $col++) {
    $newSheet.Columns.Item($col).AutoFit() | Out-Null
}
$newSheet.Range("A1:F4").Select()

$wb.Save()
